$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update A3 value
$ws.Range("A3").Value = 30

# Row 2 new values (C2:K2)
$ws.Range("C2").Value = -0.3955679064921799
$ws.Range("D2").Value = -0.3135526085529318
$ws.Range("E2").Value = -0.219670687897636
$ws.Range("F2").Value = -0.1582138509973569
$ws.Range("G2").Value = -0.1253459904459567
$ws.Range("H2").Value = -0.3248105319531392
$ws.Range("I2").Value = -0.186905187417177
$ws.Range("J2").Value = -0.1294736865881621
$ws.Range("K2").Value = -0.09242948704970184

# Row 3 new values (C3:K3)
$ws.Range("C3").Value = -0.4535728184682561
$ws.Range("D3").Value = -0.3809226417101971
$ws.Range("E3").Value = -0.3007008387427232
$ws.Range("F3").Value = -0.2529676474746937
$ws.Range("G3").Value = -0.2199545933085434
$ws.Range("H3").Value = -0.387860987107398
$ws.Range("I3").Value = -0.2678362558876907
$ws.Range("J3").Value = -0.2380870939725943
$ws.Range("K3").Value = -0.1868924887586276
